# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (E, K) from 1 to 3 for all
# data rows, and refresh the dependent NATMI-computed statistics
# (G, H, I, J, M, N, O, P, Q, R, S, T) that were recalculated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.379281
$ws.Range("H2").Value = 22.137843
$ws.Range("I2").Value = 0.2744121884499962
$ws.Range("J2").Value = 0.2744121884499961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 22.795622
$ws.Range("N2").Value = 68.386866
$ws.Range("O2").Value = 0.2266510574407626
$ws.Range("P2").Value = 0.2266510574407626
$ws.Range("Q2").Value = 168.215300307782
$ws.Range("R2").Value = 1513.937702770038
$ws.Range("S2").Value = 0.06219581268682545
$ws.Range("T2").Value = 0.06219581268682545
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.379281
$ws.Range("H3").Value = 22.137843
$ws.Range("I3").Value = 0.2744121884499962
$ws.Range("J3").Value = 0.2744121884499961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8340480000000001
$ws.Range("N3").Value = 2.502144
$ws.Range("O3").Value = 0.00829272661023916
$ws.Range("P3").Value = 0.00829272661023916
$ws.Range("Q3").Value = 6.154674559488001
$ws.Range("R3").Value = 55.39207103539201
$ws.Range("S3").Value = 0.002275625257333246
$ws.Range("T3").Value = 0.002275625257333246
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.379281
$ws.Range("H4").Value = 22.137843
$ws.Range("I4").Value = 0.2744121884499962
$ws.Range("J4").Value = 0.2744121884499961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.344218666666668
$ws.Range("N4").Value = 28.032656
$ws.Range("O4").Value = 0.09290718374597164
$ws.Range("P4").Value = 0.09290718374597164
$ws.Range("Q4").Value = 68.95361526677867
$ws.Range("R4").Value = 620.582537401008
$ws.Range("S4").Value = 0.02549486361445799
$ws.Range("T4").Value = 0.02549486361445798
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.379281
$ws.Range("H5").Value = 22.137843
$ws.Range("I5").Value = 0.2744121884499962
$ws.Range("J5").Value = 0.2744121884499961
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 67.601958
$ws.Range("N5").Value = 202.805874
$ws.Range("O5").Value = 0.6721490322030266
$ws.Range("P5").Value = 0.6721490322030266
$ws.Range("Q5").Value = 498.853844232198
$ws.Range("R5").Value = 4489.684598089782
$ws.Range("S5").Value = 0.1844458868913795
$ws.Range("T5").Value = 0.1844458868913794
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 13.29805733333333
$ws.Range("H6").Value = 39.894172
$ws.Range("I6").Value = 0.4945128143207339
$ws.Range("J6").Value = 0.4945128143207338
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 22.795622
$ws.Range("N6").Value = 68.386866
$ws.Range("O6").Value = 0.2266510574407626
$ws.Range("P6").Value = 0.2266510574407626
$ws.Range("Q6").Value = 303.1374883049946
$ws.Range("R6").Value = 2728.237394744952
$ws.Range("S6").Value = 0.1120818522838018
$ws.Range("T6").Value = 0.1120818522838018
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 13.29805733333333
$ws.Range("H7").Value = 39.894172
$ws.Range("I7").Value = 0.4945128143207339
$ws.Range("J7").Value = 0.4945128143207338
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8340480000000001
$ws.Range("N7").Value = 2.502144
$ws.Range("O7").Value = 0.00829272661023916
$ws.Range("P7").Value = 0.00829272661023916
$ws.Range("Q7").Value = 11.091218122752
$ws.Range("R7").Value = 99.82096310476801
$ws.Range("S7").Value = 0.004100859574421807
$ws.Range("T7").Value = 0.004100859574421806
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.29805733333333
$ws.Range("H8").Value = 39.894172
$ws.Range("I8").Value = 0.4945128143207339
$ws.Range("J8").Value = 0.4945128143207338
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 9.344218666666668
$ws.Range("N8").Value = 28.032656
$ws.Range("O8").Value = 0.09290718374597164
$ws.Range("P8").Value = 0.09290718374597164
$ws.Range("Q8").Value = 124.2599555645369
$ws.Range("R8").Value = 1118.339600080832
$ws.Range("S8").Value = 0.04594379290483398
$ws.Range("T8").Value = 0.04594379290483397
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.29805733333333
$ws.Range("H9").Value = 39.894172
$ws.Range("I9").Value = 0.4945128143207339
$ws.Range("J9").Value = 0.4945128143207338
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 67.601958
$ws.Range("N9").Value = 202.805874
$ws.Range("O9").Value = 0.6721490322030266
$ws.Range("P9").Value = 0.6721490322030266
$ws.Range("Q9").Value = 898.9747133295919
$ws.Range("R9").Value = 8090.772419966327
$ws.Range("S9").Value = 0.3323863095576763
$ws.Range("T9").Value = 0.3323863095576762
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.673314666666667
$ws.Range("H10").Value = 5.019944000000001
$ws.Range("I10").Value = 0.06222529534320158
$ws.Range("J10").Value = 0.06222529534320156
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 22.795622
$ws.Range("N10").Value = 68.386866
$ws.Range("O10").Value = 0.2266510574407626
$ws.Range("P10").Value = 0.2266510574407626
$ws.Range("Q10").Value = 38.14424862838933
$ws.Range("R10").Value = 343.298237655504
$ws.Range("S10").Value = 0.0141034289891004
$ws.Range("T10").Value = 0.0141034289891004
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.673314666666667
$ws.Range("H11").Value = 5.019944000000001
$ws.Range("I11").Value = 0.06222529534320158
$ws.Range("J11").Value = 0.06222529534320156
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.8340480000000001
$ws.Range("N11").Value = 2.502144
$ws.Range("O11").Value = 0.00829272661023916
$ws.Range("P11").Value = 0.00829272661023916
$ws.Range("Q11").Value = 1.395624751104001
$ws.Range("R11").Value = 12.560622759936
$ws.Range("S11").Value = 0.0005160173625225586
$ws.Range("T11").Value = 0.0005160173625225585
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.673314666666667
$ws.Range("H12").Value = 5.019944000000001
$ws.Range("I12").Value = 0.06222529534320158
$ws.Range("J12").Value = 0.06222529534320156
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 9.344218666666668
$ws.Range("N12").Value = 28.032656
$ws.Range("O12").Value = 0.09290718374597164
$ws.Range("P12").Value = 0.09290718374597164
$ws.Range("Q12").Value = 15.63581814347378
$ws.Range("R12").Value = 140.722363291264
$ws.Range("S12").Value = 0.005781176948098182
$ws.Range("T12").Value = 0.005781176948098181
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.673314666666667
$ws.Range("H13").Value = 5.019944000000001
$ws.Range("I13").Value = 0.06222529534320158
$ws.Range("J13").Value = 0.06222529534320156
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 67.601958
$ws.Range("N13").Value = 202.805874
$ws.Range("O13").Value = 0.6721490322030266
$ws.Range("P13").Value = 0.6721490322030266
$ws.Range("Q13").Value = 113.119347816784
$ws.Range("R13").Value = 1018.074130351056
$ws.Range("S13").Value = 0.04182467204348043
$ws.Range("T13").Value = 0.04182467204348043
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.540576
$ws.Range("H14").Value = 13.621728
$ws.Range("I14").Value = 0.1688497018860685
$ws.Range("J14").Value = 0.1688497018860685
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 22.795622
$ws.Range("N14").Value = 68.386866
$ws.Range("O14").Value = 0.2266510574407626
$ws.Range("P14").Value = 0.2266510574407626
$ws.Range("Q14").Value = 103.505254158272
$ws.Range("R14").Value = 931.5472874244479
$ws.Range("S14").Value = 0.03826996348103496
$ws.Range("T14").Value = 0.03826996348103496
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.540576
$ws.Range("H15").Value = 13.621728
$ws.Range("I15").Value = 0.1688497018860685
$ws.Range("J15").Value = 0.1688497018860685
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.8340480000000001
$ws.Range("N15").Value = 2.502144
$ws.Range("O15").Value = 0.00829272661023916
$ws.Range("P15").Value = 0.00829272661023916
$ws.Range("Q15").Value = 3.787058331648
$ws.Range("R15").Value = 34.083524984832
$ws.Range("S15").Value = 0.00140022441596155
$ws.Range("T15").Value = 0.001400224415961549
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.540576
$ws.Range("H16").Value = 13.621728
$ws.Range("I16").Value = 0.1688497018860685
$ws.Range("J16").Value = 0.1688497018860685
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 9.344218666666668
$ws.Range("N16").Value = 28.032656
$ws.Range("O16").Value = 0.09290718374597164
$ws.Range("P16").Value = 0.09290718374597164
$ws.Range("Q16").Value = 42.42813501661867
$ws.Range("R16").Value = 381.853215149568
$ws.Range("S16").Value = 0.0156873502785815
$ws.Range("T16").Value = 0.0156873502785815
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.540576
$ws.Range("H17").Value = 13.621728
$ws.Range("I17").Value = 0.1688497018860685
$ws.Range("J17").Value = 0.1688497018860685
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 67.601958
$ws.Range("N17").Value = 202.805874
$ws.Range("O17").Value = 0.6721490322030266
$ws.Range("P17").Value = 0.6721490322030266
$ws.Range("Q17").Value = 306.951828047808
$ws.Range("R17").Value = 2762.566452430272
$ws.Range("S17").Value = 0.1134921637104905
$ws.Range("T17").Value = 0.1134921637104905
